# tests/assets/excel_data.xlsx - "done what Rouven said"
#
# Patient P004's Sex/Living answers are retracted (cleared), and a new
# "Age of onset" value is recorded for patient P002's second condition
# (increased size of head) on the "conditions" sheet. Finally, the user
# ends up with the "basic info" sheet active with C5 selected (having
# last looked at F12 on "conditions").

$wb = $excel.ActiveWorkbook

$wsBasic = $wb.Worksheets.Item("basic info")
$wsConditions = $wb.Worksheets.Item("conditions")

# P004's Sex ("m") and Living ("Yes") answers are removed.
$wsBasic.Range("B5").ClearContents()
$wsBasic.Range("C5").ClearContents()

# New Age-of-onset value for P002 / "increased size of head".
$wsConditions.Range("E4").Value = "P7Y5M8D"

# Leave a selection behind on "conditions" ...
$wsConditions.Range("F12").Select()

# ... then switch to "basic info" and select C5, which is where the
# file was left (and saved) from.
$wsBasic.Activate()
$wsBasic.Range("C5").Select()
